$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '42.467.40'
$ws.Range("E2").Value = '  +1.76%  '
$ws.Range("D3").Value = '2.291.27'
$ws.Range("E3").Value = '  +1.18%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '157.15'
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  +15,611.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '307.57'
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = '  +1.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '95.87'
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = '  +4.83%  '
$ws.Range("E8").Value = '  +0.30%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.496'
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = '  +3.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '36.12'
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = '  +12.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0804'
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = '  +1.14%  '
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.74'
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = '  +2.51%  '
$ws.Range("D15").Value = '2.643.51'
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.54'
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = '  +2.55%  '
$ws.Range("D17").Value = '2.298.66'
$ws.Range("E17").Value = '  +1.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.801'
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = '  +5.44%  '
$ws.Range("D19").Value = '42.366.93'
$ws.Range("E19").Value = '  +1.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.72'
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = '  +2.62%  '
$ws.Range("D21").Value = '0.0₃0918'
$ws.Range("E21").Value = '  +1.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.02'
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  +2.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.05'
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '243.36'
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = '  +1.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.61'
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.95'
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = '  +2.34%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.11'
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35.95'
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = '  +3.68%  '
$ws.Range("E30").Value = '  +0.98%  '
$ws.Range("E31").Value = '  +1.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '161.51'
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("E33").Value = '  +3.84%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0756'
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = '  +1.68%  '
$ws.Range("E36").Value = '  +3.03%  '
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.27'
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = '  +3.99%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.108'
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = '  +4.05%  '
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.85'
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = '  +3.28%  '
$ws.Range("E41").Value = '  -0.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.18'
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = '  +6.87%  '
$ws.Range("D43").Value = '2.010.57'
$ws.Range("E43").Value = '  -2.47%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.29'
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = '  +12.04%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.54'
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = '  +0.84%  '
$ws.Range("E46").Value = '  +2.86%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.00'
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = '  +4.90%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.14'
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.88'
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = '  +4.32%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.55'
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = '  +2.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.00'
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = '  +0.33%  '
